# Add three new MIxS columns (Collector Sample ID, Description Of Collection
# Method, Taxonomy ID) to the "Sample" sheet, as optional columns D, E, F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# The sheet is protected; unprotect while we edit, re-protect at the end.
$ws.Unprotect()

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("D1").Value = "Collector Sample ID (optional)"
$ws.Range("E1").Value = "Description Of Collection Method (optional)"
$ws.Range("F1").Value = "Taxonomy ID (optional)"

# --- Description row (row 2) ----------------------------------------------
$ws.Range("D2").Value = 'Unique name assigned to the sample by the COLLECTOR or COLLECTOR_AFFILIATION. Do not use spaces or special characters, other than hyphens and underscores ("-" and "_") i.e do not use #, !,^*, etc. Free text.'
$ws.Range("E2").Value = "A detailed as possible description of the sample collection methods"
$ws.Range("F2").Value = 'A valid NCBI TAXON_ID to the species level ismandatory in order to submit data to public repositories. The species name in themanifest must be identical to that listed in the "current name" box in the T axonomy Browser for that species.'

# --- Example row (row 3) ----------------------------------------------------
$ws.Range("D3").Value = "e.g. UDUK0000331"
$ws.Range("E3").Value = "e.g. Caught with fiber net within densely wooded area, and immediately placed into the collection container"
$ws.Range("F3").Value = "e.g. 458489"

# --- Copy header/description/example formatting from column C -------------
$ws.Range("C1:C3").Copy()
$ws.Range("D1:D3").PasteSpecial(-4122)
$ws.Range("C1:C3").Copy()
$ws.Range("E1:E3").PasteSpecial(-4122)
$ws.Range("C1:C3").Copy()
$ws.Range("F1:F3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Extend the "FILL OUT INFORMATION BELOW THIS ROW" merged banner -------
$ws.Range("A4:C4").UnMerge()
$ws.Range("A4:F4").Merge()

# --- Column widths (best-fit to content, matching columns A-C style) ------
$ws.Columns.Item(4).ColumnWidth = 184.6
$ws.Columns.Item(5).ColumnWidth = 95.75
$ws.Columns.Item(6).ColumnWidth = 207.92

# --- Extend conditional formatting (error highlighting) over new columns --
$fcs = $ws.Range("A1:C4").FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A1:F4"))

# --- Data validation for the new columns -----------------------------------
$dFormula = 'AND(SUM(--ISNUMBER(FIND(MID(D5, ROW(INDIRECT("1:"&LEN(D5))), 1), "abcdefghijklmnopqrstuvwxyz"))) > 0)'
$ws.Range("D5:D1005").Validation.Add(7, 1, 1, $dFormula)

$eFormula = 'AND(SUM(--ISNUMBER(FIND(MID(E5, ROW(INDIRECT("1:"&LEN(E5))), 1), "abcdefghijklmnopqrstuvwxyz"))) > 0)'
$ws.Range("E5:E1005").Validation.Add(7, 1, 1, $eFormula)

$fFormula = 'AND(ISNUMBER(F5+0), INT(F5+0)=F5+0)'
$ws.Range("F5:F1005").Validation.Add(7, 1, 1, $fFormula)

# --- Re-protect the sheet, same as before ----------------------------------
$ws.Protect()
